$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.07271233333333334
$ws.Range("H2").Value = 0.218137
$ws.Range("I2").Value = 0.004171225362010892
$ws.Range("J2").Value = 0.004171225362010893
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 0.8105131819673335
$ws.Range("R2").Value = 7.294618637706002
$ws.Range("S2").Value = 0.001082352095405582
$ws.Range("T2").Value = 0.001082352095405582
$ws.Range("G3").Value = 0.07271233333333334
$ws.Range("H3").Value = 0.218137
$ws.Range("I3").Value = 0.004171225362010892
$ws.Range("J3").Value = 0.004171225362010893
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 2.013077086363
$ws.Range("R3").Value = 18.117693777267
$ws.Range("S3").Value = 0.00268824523908332
$ws.Range("T3").Value = 0.002688245239083321
$ws.Range("G4").Value = 0.07271233333333334
$ws.Range("H4").Value = 0.218137
$ws.Range("I4").Value = 0.004171225362010892
$ws.Range("J4").Value = 0.004171225362010893
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 0.3000080091778889
$ws.Range("R4").Value = 2.700072082601001
$ws.Range("S4").Value = 0.0004006280275219907
$ws.Range("T4").Value = 0.0004006280275219908
$ws.Range("I5").Value = 0.5387060579248023
$ws.Range("J5").Value = 0.5387060579248023
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 104.6762817301287
$ws.Range("R5").Value = 942.086535571158
$ws.Range("S5").Value = 0.1397837757491723
$ws.Range("T5").Value = 0.1397837757491723
$ws.Range("I6").Value = 0.5387060579248023
$ws.Range("J6").Value = 0.5387060579248023
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.3471819117400906
$ws.Range("T6").Value = 0.3471819117400906
$ws.Range("I7").Value = 0.5387060579248023
$ws.Range("J7").Value = 0.5387060579248023
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.05174037043553947
$ws.Range("T7").Value = 0.05174037043553947
$ws.Range("I8").Value = 0.4571227167131868
$ws.Range("J8").Value = 0.4571227167131868
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 88.82377611315202
$ws.Range("R8").Value = 799.4139850183682
$ws.Range("S8").Value = 0.1186144807226357
$ws.Range("T8").Value = 0.1186144807226357
$ws.Range("I9").Value = 0.4571227167131868
$ws.Range("J9").Value = 0.4571227167131868
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.2946035901279238
$ws.Range("T9").Value = 0.2946035901279238
$ws.Range("I10").Value = 0.4571227167131868
$ws.Range("J10").Value = 0.4571227167131868
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.04390464586262734
$ws.Range("T10").Value = 0.04390464586262734

Write-Output "Updated 76 cells"
